# Update the NBA player roster table (A2:C17) on Sheet1 to the new lineup.
# Two players were removed (Nikola Jovic, Klay Thompson), two were added
# (Gradey Dick, Jordan Clarkson), and the remaining rows were reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Tobias Harris",           "SF,PF", "Detroit Pistons"),
    @("Kyle Kuzma",               "PF",    "Milwaukee Bucks"),
    @("Gradey Dick",              "SG,SF", "Toronto Raptors"),
    @("Joel Embiid",              "C",     "Philadelphia 76ers"),
    @("Jordan Poole",             "PG,SG", "Washington Wizards"),
    @("Shai Gilgeous-Alexander",  "PG,SG", "Oklahoma City Thunder"),
    @("Kyrie Irving",             "PG,SG", "Dallas Mavericks"),
    @("CJ McCollum",              "PG,SG", "New Orleans Pelicans"),
    @("Lauri Markkanen",          "SF,PF", "Utah Jazz"),
    @("Zach LaVine",              "SG,SF", "Sacramento Kings"),
    @("Toumani Camara",           "SF,PF", "Portland Trail Blazers"),
    @("Rui Hachimura",            "SF,PF", "Los Angeles Lakers"),
    @("Bam Adebayo",              "PF,C",  "Miami Heat"),
    @("Jordan Clarkson",          "SG,SF", "Utah Jazz"),
    @("Jamal Murray",             "PG,SG", "Denver Nuggets"),
    @("John Collins",             "PF,C",  "Utah Jazz")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
